$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "35.097.98"
$ws.Cells.Item(2, 5).Value = "  -0.89%  "

$ws.Cells.Item(3, 4).Value = "1.895.72"
$ws.Cells.Item(3, 5).Value = "  -0.47%  "

$ws.Cells.Item(4, 5).Value = "  -0.40%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "252.37"
$ws.Cells.Item(5, 5).Value = "  +2.28%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.697"
$ws.Cells.Item(6, 5).Value = "  -0.41%  "

$ws.Cells.Item(7, 5).Value = "  -0.39%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "41.48"
$ws.Cells.Item(8, 5).Value = "  +1.99%  "

$ws.Cells.Item(9, 5).Value = "  +0.78%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0753"
$ws.Cells.Item(10, 5).Value = "  +3.74%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0978"
$ws.Cells.Item(11, 5).Value = "  -1.18%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "13.09"
$ws.Cells.Item(12, 5).Value = "  +3.75%  "

$ws.Cells.Item(13, 4).Value = "2.172.78"
$ws.Cells.Item(13, 5).Value = "  -0.43%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.726"
$ws.Cells.Item(14, 5).Value = "  +1.82%  "

$ws.Cells.Item(15, 5).Value = "  +1.53%  "

$ws.Cells.Item(16, 4).Value = "1.906.37"
$ws.Cells.Item(16, 5).Value = "  -0.55%  "

$ws.Cells.Item(17, 4).Value = "35.115.13"
$ws.Cells.Item(17, 5).Value = "  -0.78%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "74.12"
$ws.Cells.Item(18, 5).Value = "  +1.50%  "

$ws.Cells.Item(19, 4).Value = "0.0₃0835"
$ws.Cells.Item(19, 5).Value = "  +1.23%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "252.07"
$ws.Cells.Item(20, 5).Value = "  +4.34%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "12.94"
$ws.Cells.Item(21, 5).Value = "  +0.09%  "

$ws.Cells.Item(22, 5).Value = "  -1.50%  "

$ws.Cells.Item(23, 5).Value = "  -0.41%  "

$ws.Cells.Item(24, 5).Value = "  +4.58%  "

$ws.Cells.Item(25, 5).Value = "  -2.12%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "168.39"
$ws.Cells.Item(26, 5).Value = "  -0.45%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "8.55"
$ws.Cells.Item(27, 5).Value = "  -1.19%  "

$ws.Cells.Item(28, 5).Value = "  -2.37%  "

$ws.Cells.Item(29, 5).Value = "  -1.74%  "

$ws.Cells.Item(30, 4).Value = "4.128.32"
$ws.Cells.Item(30, 5).Value = "  -0.34%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.31"
$ws.Cells.Item(31, 5).Value = "  +1.83%  "

$ws.Cells.Item(32, 2).Value = "Hedera"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.0593"
$ws.Cells.Item(32, 5).Value = "  +3.37%  "

$ws.Cells.Item(33, 2).Value = "TrustWalletToken"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.62"
$ws.Cells.Item(33, 5).Value = "  +10.11%  "

$ws.Cells.Item(34, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.23"
$ws.Cells.Item(34, 5).Value = "  +1.26%  "

$ws.Cells.Item(35, 2).Value = "WEMIXToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.87"
$ws.Cells.Item(35, 5).Value = "  -0.36%  "

$ws.Cells.Item(36, 5).Value = "  -0.46%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.848"
$ws.Cells.Item(37, 5).Value = "  -7.56%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.01"
$ws.Cells.Item(38, 5).Value = "  -0.73%  "

$ws.Cells.Item(39, 5).Value = "  +5.53%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "98.74"
$ws.Cells.Item(40, 5).Value = "  +1.79%  "

$ws.Cells.Item(41, 5).Value = "  +2.43%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0659"
$ws.Cells.Item(42, 5).Value = "  +0.41%  "

$ws.Cells.Item(43, 5).Value = "  -0.63%  "

$ws.Cells.Item(44, 2).Value = "Maker"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(44, 4).Value = "1.301.74"
$ws.Cells.Item(44, 5).Value = "  -4.08%  "

$ws.Cells.Item(45, 2).Value = "RenderToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.38"
$ws.Cells.Item(45, 5).Value = "  -1.37%  "

$ws.Cells.Item(46, 5).Value = "  -0.08%  "

$ws.Cells.Item(47, 5).Value = "  -1.57%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "6.57"
$ws.Cells.Item(48, 5).Value = "  +0.55%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0769"
$ws.Cells.Item(49, 5).Value = "  +8.35%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "11.93"
$ws.Cells.Item(50, 5).Value = "  -3.12%  "

$ws.Cells.Item(51, 5).Value = "  -7.51%  "
